$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Acelga" (Macroferia Regional de Talca) is
# inserted as row 338; every existing row from 338 onward shifts down by one
# (dimension grows from A1:R478 to A1:R479).
$ws.Rows.Item(338).Insert()

$ws.Cells.Item(338, 1).Value = 5
$ws.Cells.Item(338, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(338, 3).Value = 'Maule'
$ws.Cells.Item(338, 4).Value = 45141
$ws.Cells.Item(338, 5).Value = 7
$ws.Cells.Item(338, 6).Value = 100112009
$ws.Cells.Item(338, 7).Value = 'Acelga'
$ws.Cells.Item(338, 8).Value = 'Sin especificar'
$ws.Cells.Item(338, 9).Value = 'Primera'
$ws.Cells.Item(338, 10).Value = 500
$ws.Cells.Item(338, 11).Value = 1500
$ws.Cells.Item(338, 12).Value = 1500
$ws.Cells.Item(338, 13).Value = 1500
$ws.Cells.Item(338, 14).Value = '$/docena de atados (4 kilos)'
$ws.Cells.Item(338, 15).Value = 'Región del Maule'
$ws.Cells.Item(338, 16).Value = 375
$ws.Cells.Item(338, 17).Value = 4
$ws.Cells.Item(338, 18).Value = 'Hortaliza'
